# koboforms/section9.xlsx: section8 "Savings Discussion" -> section9 "YourGumGame"
# Rewrites the survey questions (q1,q2,q4,q5 reworded; new q6 added; q7 reworded;
# former select_one-yes_no row becomes a plain text question), adds an
# "appearance" column with "multiline" for every question, updates the
# settings sheet (form_title/form_id), and moves the active tab/selection
# from "settings" to "survey".

$wb = $excel.ActiveWorkbook

# ---- settings sheet: form_title / form_id ----
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "YourGumGame"
$settings.Range("B2").Value = "section9"

# ---- survey sheet: questions + new "appearance" column ----
$survey = $wb.Worksheets.Item("survey")

# Clone formatting for the new cells from their row/column neighbours
# before filling in the new values, so the new column matches the existing
# header/body look (bold Arial header row, plain Arial body rows).
$survey.Range("D1").Copy()
$survey.Range("E1").PasteSpecial(-4122)

$survey.Range("D3").Copy()
$survey.Range("E3:E8").PasteSpecial(-4122)

$survey.Range("A3").Copy()
$survey.Range("A7").PasteSpecial(-4122)

# New header for column E
$survey.Range("E1").Value = "appearance"

# Row 2 (email) is untouched aside from gaining no "appearance" value.

# Row 3: q1
$survey.Range("B3").Value = "q1"
$survey.Range("C3").Value = "Describe the scenario for the game"
$survey.Range("E3").Value = "multiline"

# Row 4: q2
$survey.Range("B4").Value = "q2"
$survey.Range("C4").Value = "Where do your parameters come from?"
$survey.Range("E4").Value = "multiline"

# Row 5: q4
$survey.Range("B5").Value = "q4"
$survey.Range("C5").Value = "Where do you make reasonable assumptions instead of getting parameters from references?"
$survey.Range("E5").Value = "multiline"

# Row 6: q5
$survey.Range("B6").Value = "q5"
$survey.Range("C6").Value = "What are the general instructions for the game?"
$survey.Range("E6").Value = "multiline"

# Row 7: was "select_one yes_no" / partial_useful -> now plain text q6
$survey.Range("A7").Value = "text"
$survey.Range("B7").Value = "q6"
$survey.Range("C7").Value = "Describe the investment that creates value, and how"
$survey.Range("E7").Value = "multiline"

# Row 8: q7
$survey.Range("B8").Value = "q7"
$survey.Range("C8").Value = "Describe how since insurance simply moves money between years, for a fee, and how someone would tell if it is worthwhile or not."
$survey.Range("E8").Value = "multiline"

# ---- view state: make "survey" the active/selected tab, settings loses it ----
$survey.Activate()
$survey.Range("C7").Select()
